# Permission Processing On Api
# Remove the four obsolete attachment-related permission rows (rows 12-15,
# which referenced the "上传/下载/我的文件/读取文件" permIds) from the
# DATA-PERM sheet. Deleting the entire rows shifts the trailing S_VIEW
# block (old rows 19-21) up to rows 15-17 and lets the host re-pack the
# shared-string table automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 12:15 (1-based) hold the permId entries for 上传/下载/我的文件/读取文件.
$ws.Range("A12:K15").EntireRow.Delete()

# Re-apply the on-screen selection that Excel leaves behind after deleting
# those rows (matches the post-edit view state).
$ws.Range("A12:XFD15").Select()

# Best-effort: keep the window position in sync with the authored view.
$excel.ActiveWindow.Left = 43680
$excel.ActiveWindow.Top = -10500
